$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 8-10 (ECs as a sending cluster is removed entirely)
$ws.Rows("8:10").Delete()

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Pdgfc"
$ws.Range("C2").Value = "Pdgfrb"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.441874
$ws.Range("H2").Value = 7.325622
$ws.Range("I2").Value = 0.7339587032246254
$ws.Range("J2").Value = 0.7339587032246254
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 6.07605
$ws.Range("N2").Value = 18.22815
$ws.Range("O2").Value = 0.0302610603580868
$ws.Range("P2").Value = 0.0302610603580868
$ws.Range("Q2").Value = 14.8369485177
$ws.Range("R2").Value = 133.5325366593
$ws.Range("S2").Value = 0.02221036861862351
$ws.Range("T2").Value = 0.0222103686186235

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Pdgfc"
$ws.Range("C3").Value = "Pdgfrb"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.441874
$ws.Range("H3").Value = 7.325622
$ws.Range("I3").Value = 0.7339587032246254
$ws.Range("J3").Value = 0.7339587032246254
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 60.13240533333334
$ws.Range("N3").Value = 180.397216
$ws.Range("O3").Value = 0.2994824511432495
$ws.Range("P3").Value = 0.2994824511432494
$ws.Range("Q3").Value = 146.835757140928
$ws.Range("R3").Value = 1321.521814268352
$ws.Range("S3").Value = 0.2198077514796316
$ws.Range("T3").Value = 0.2198077514796316

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Pdgfc"
$ws.Range("C4").Value = "Pdgfrb"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.441874
$ws.Range("H4").Value = 7.325622
$ws.Range("I4").Value = 0.7339587032246254
$ws.Range("J4").Value = 0.7339587032246254
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 134.5792873333333
$ws.Range("N4").Value = 403.737862
$ws.Range("O4").Value = 0.6702564884986638
$ws.Range("P4").Value = 0.6702564884986637
$ws.Range("Q4").Value = 328.625662677796
$ws.Range("R4").Value = 2957.630964100164
$ws.Range("S4").Value = 0.4919405831263703
$ws.Range("T4").Value = 0.4919405831263703

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Pdgfc"
$ws.Range("C5").Value = "Pdgfrb"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.885117
$ws.Range("H5").Value = 2.655351
$ws.Range("I5").Value = 0.2660412967753745
$ws.Range("J5").Value = 0.2660412967753745
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 6.07605
$ws.Range("N5").Value = 18.22815
$ws.Range("O5").Value = 0.0302610603580868
$ws.Range("P5").Value = 0.0302610603580868
$ws.Range("Q5").Value = 5.37801514785
$ws.Range("R5").Value = 48.40213633065
$ws.Range("S5").Value = 0.008050691739463291
$ws.Range("T5").Value = 0.00805069173946329

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Pdgfc"
$ws.Range("C6").Value = "Pdgfrb"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.885117
$ws.Range("H6").Value = 2.655351
$ws.Range("I6").Value = 0.2660412967753745
$ws.Range("J6").Value = 0.2660412967753745
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 60.13240533333334
$ws.Range("N6").Value = 180.397216
$ws.Range("O6").Value = 0.2994824511432495
$ws.Range("P6").Value = 0.2994824511432494
$ws.Range("Q6").Value = 53.22421421142401
$ws.Range("R6").Value = 479.017927902816
$ws.Range("S6").Value = 0.07967469966361783
$ws.Range("T6").Value = 0.0796746996636178

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Pdgfc"
$ws.Range("C7").Value = "Pdgfrb"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.885117
$ws.Range("H7").Value = 2.655351
$ws.Range("I7").Value = 0.2660412967753745
$ws.Range("J7").Value = 0.2660412967753745
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 134.5792873333333
$ws.Range("N7").Value = 403.737862
$ws.Range("O7").Value = 0.6702564884986638
$ws.Range("P7").Value = 0.6702564884986637
$ws.Range("Q7").Value = 119.118415066618
$ws.Range("R7").Value = 1072.065735599562
$ws.Range("S7").Value = 0.1783159053722934
$ws.Range("T7").Value = 0.1783159053722934
